$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "LEGO"
$ws.Range("C4").Value = "Cerru"
$ws.Range("D4").Value = 12000
$ws.Range("E4").Value = 10
